$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.06"
$ws.Range("E2").Value = "'0.88%"
$ws.Range("D3").Value = "'41.28"
$ws.Range("E3").Value = "'2.75%"
$ws.Range("D4").Value = "'5.706"
$ws.Range("E4").Value = "'-1.93%"
$ws.Range("D5").Value = "'0.08089"
$ws.Range("E5").Value = "'0.82%"
$ws.Range("D6").Value = "'2.040"
$ws.Range("E6").Value = "'4.95%"
$ws.Range("D7").Value = "'8.726"
$ws.Range("E7").Value = "'0.20%"
$ws.Range("D8").Value = "'4.508"
$ws.Range("E8").Value = "'-1.47%"
$ws.Range("E9").Value = "'-0.72%"
$ws.Range("D11").Value = "'0.1252"
$ws.Range("E11").Value = "'0.09%"
$ws.Range("D12").Value = "'0.1945"
$ws.Range("E12").Value = "'-0.79%"
$ws.Range("D13").Value = "'8.323"
$ws.Range("E13").Value = "'-6.76%"
$ws.Range("D14").Value = "'0.09318"
$ws.Range("E14").Value = "'1.21%"
$ws.Range("D15").Value = "'0.03664"
$ws.Range("D16").Value = "'0.1053"
$ws.Range("E16").Value = "'9.45%"
$ws.Range("D17").Value = "'0.001307"
$ws.Range("E17").Value = "'0.68%"
$ws.Range("D18").Value = "'0.006137"
$ws.Range("E18").Value = "'-1.24%"
$ws.Range("E19").Value = "'0.32%"
$ws.Range("E21").Value = "'0.77%"
$ws.Range("D22").Value = "'0.2648"
$ws.Range("E22").Value = "'9.44%"
$ws.Range("D23").Value = "'0.04433"
$ws.Range("E23").Value = "'0.64%"
$ws.Range("D24").Value = "'0.001259"
$ws.Range("E24").Value = "'-0.10%"
$ws.Range("D25").Value = "'0.004324"
$ws.Range("E25").Value = "'-0.67%"
$ws.Range("E26").Value = "'8.38%"
$ws.Range("D39").Value = "'0.02816"
$ws.Range("E39").Value = "'16.39%"
$ws.Range("D40").Value = "'0.05471"
$ws.Range("D41").Value = "'0.007576"
$ws.Range("E41").Value = "'1.21%"
$ws.Range("D42").Value = "'0.009953"
$ws.Range("E42").Value = "'15.60%"
$ws.Range("D43").Value = "'0.1421"
$ws.Range("E43").Value = "'0.20%"
$ws.Range("D44").Value = "'0.002132"
$ws.Range("E44").Value = "'1.09%"
$ws.Range("D45").Value = "'0.01186"
$ws.Range("E45").Value = "'25.11%"
$ws.Range("D46").Value = "'0.00006750"
$ws.Range("E46").Value = "'-2.07%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.43%"
$ws.Range("D48").Value = "'0.003077"
$ws.Range("E48").Value = "'-2.58%"
$ws.Range("D49").Value = "'0.002278"
$ws.Range("E49").Value = "'59.79%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.43%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.43%"
